$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet 1: Recommandations — refresh rows 2-44 with latest BRVM data ---
# Row 2: BRVM - SERVICES PUBLICS
$ws1.Cells.Item(2,1).Value2 = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2,2).Value2 = 0
$ws1.Cells.Item(2,3).Value2 = 6
$ws1.Cells.Item(2,4).Value2 = 2417.55
$ws1.Cells.Item(2,5).Value2 = 102.31
$ws1.Cells.Item(2,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(2,7).Value2 = "➖ Neutre"

# Row 3: CFAO MOTORS CI
$ws1.Cells.Item(3,1).Value2 = "CFAO MOTORS CI"
$ws1.Cells.Item(3,2).Value2 = 0
$ws1.Cells.Item(3,3).Value2 = 3
$ws1.Cells.Item(3,4).Value2 = 1945
$ws1.Cells.Item(3,5).Value2 = 645
$ws1.Cells.Item(3,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(3,7).Value2 = "➖ Neutre"

# Row 4: BRVM - AUTRES SECTEURS
$ws1.Cells.Item(4,1).Value2 = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(4,2).Value2 = 0
$ws1.Cells.Item(4,3).Value2 = 3
$ws1.Cells.Item(4,4).Value2 = 1907.97
$ws1.Cells.Item(4,5).Value2 = 641.5700000000001
$ws1.Cells.Item(4,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(4,7).Value2 = "➖ Neutre"

# Row 5: SAFCA CI
$ws1.Cells.Item(5,1).Value2 = "SAFCA CI"
$ws1.Cells.Item(5,2).Value2 = 0
$ws1.Cells.Item(5,3).Value2 = 2
$ws1.Cells.Item(5,4).Value2 = 1905
$ws1.Cells.Item(5,5).Value2 = 985
$ws1.Cells.Item(5,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(5,7).Value2 = "➖ Neutre"

# Row 6: NEI-CEDA CI
$ws1.Cells.Item(6,1).Value2 = "NEI-CEDA CI"
$ws1.Cells.Item(6,2).Value2 = 0
$ws1.Cells.Item(6,3).Value2 = 3
$ws1.Cells.Item(6,4).Value2 = 1780
$ws1.Cells.Item(6,5).Value2 = 595
$ws1.Cells.Item(6,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(6,7).Value2 = "➖ Neutre"

# Row 7: SETAO CI
$ws1.Cells.Item(7,1).Value2 = "SETAO CI"
$ws1.Cells.Item(7,2).Value2 = 0
$ws1.Cells.Item(7,3).Value2 = 3
$ws1.Cells.Item(7,4).Value2 = 1750
$ws1.Cells.Item(7,5).Value2 = 580
$ws1.Cells.Item(7,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(7,7).Value2 = "➖ Neutre"

# Row 8: UNIWAX CI
$ws1.Cells.Item(8,1).Value2 = "UNIWAX CI"
$ws1.Cells.Item(8,2).Value2 = 0
$ws1.Cells.Item(8,3).Value2 = 3
$ws1.Cells.Item(8,4).Value2 = 1725
$ws1.Cells.Item(8,5).Value2 = 580
$ws1.Cells.Item(8,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(8,7).Value2 = "➖ Neutre"

# Row 9: AIR LIQUIDE CI
$ws1.Cells.Item(9,1).Value2 = "AIR LIQUIDE CI"
$ws1.Cells.Item(9,2).Value2 = 0
$ws1.Cells.Item(9,3).Value2 = 3
$ws1.Cells.Item(9,4).Value2 = 1590
$ws1.Cells.Item(9,5).Value2 = 545
$ws1.Cells.Item(9,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(9,7).Value2 = "➖ Neutre"

# Row 10: BRVM - DISTRIBUTION
$ws1.Cells.Item(10,1).Value2 = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(10,2).Value2 = 0
$ws1.Cells.Item(10,3).Value2 = 3
$ws1.Cells.Item(10,4).Value2 = 1089.18
$ws1.Cells.Item(10,5).Value2 = 364.57
$ws1.Cells.Item(10,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(10,7).Value2 = "➖ Neutre"

# Row 11: BRVM - TRANSPORT
$ws1.Cells.Item(11,1).Value2 = "BRVM - TRANSPORT"
$ws1.Cells.Item(11,2).Value2 = 0
$ws1.Cells.Item(11,3).Value2 = 3
$ws1.Cells.Item(11,4).Value2 = 1050.07
$ws1.Cells.Item(11,5).Value2 = 352.47
$ws1.Cells.Item(11,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(11,7).Value2 = "➖ Neutre"

# Row 12: BRVM - AGRICULTURE
$ws1.Cells.Item(12,1).Value2 = "BRVM - AGRICULTURE"
$ws1.Cells.Item(12,2).Value2 = 0
$ws1.Cells.Item(12,3).Value2 = 3
$ws1.Cells.Item(12,4).Value2 = 955.28
$ws1.Cells.Item(12,5).Value2 = 321.44
$ws1.Cells.Item(12,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(12,7).Value2 = "➖ Neutre"

# Row 13: BRVM - INDUSTRIE
$ws1.Cells.Item(13,1).Value2 = "BRVM - INDUSTRIE"
$ws1.Cells.Item(13,2).Value2 = 0
$ws1.Cells.Item(13,3).Value2 = 3
$ws1.Cells.Item(13,4).Value2 = 786.15
$ws1.Cells.Item(13,5).Value2 = 262.61
$ws1.Cells.Item(13,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(13,7).Value2 = "➖ Neutre"

# Row 14: BRVM - CONSOMMATION DE BASE
$ws1.Cells.Item(14,1).Value2 = "BRVM - CONSOMMATION DE BASE"
$ws1.Cells.Item(14,2).Value2 = 0
$ws1.Cells.Item(14,3).Value2 = 3
$ws1.Cells.Item(14,4).Value2 = 648.08
$ws1.Cells.Item(14,5).Value2 = 216.86
$ws1.Cells.Item(14,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(14,7).Value2 = "➖ Neutre"

# Row 15: BRVM-PRINCIPAL
$ws1.Cells.Item(15,1).Value2 = "BRVM-PRINCIPAL"
$ws1.Cells.Item(15,2).Value2 = 0
$ws1.Cells.Item(15,3).Value2 = 3
$ws1.Cells.Item(15,4).Value2 = 563.9299999999999
$ws1.Cells.Item(15,5).Value2 = 188.77
$ws1.Cells.Item(15,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(15,7).Value2 = "➖ Neutre"

# Row 16: BRVM - INDUSTRIELS
$ws1.Cells.Item(16,1).Value2 = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(16,2).Value2 = 0
$ws1.Cells.Item(16,3).Value2 = 3
$ws1.Cells.Item(16,4).Value2 = 413.66
$ws1.Cells.Item(16,5).Value2 = 138.37
$ws1.Cells.Item(16,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(16,7).Value2 = "➖ Neutre"

# Row 17: BRVM-PRESTIGE
$ws1.Cells.Item(17,1).Value2 = "BRVM-PRESTIGE"
$ws1.Cells.Item(17,2).Value2 = 0
$ws1.Cells.Item(17,3).Value2 = 3
$ws1.Cells.Item(17,4).Value2 = 386.14
$ws1.Cells.Item(17,5).Value2 = 128.94
$ws1.Cells.Item(17,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(17,7).Value2 = "➖ Neutre"

# Row 18: BRVM - FINANCES
$ws1.Cells.Item(18,1).Value2 = "BRVM - FINANCES"
$ws1.Cells.Item(18,2).Value2 = 0
$ws1.Cells.Item(18,3).Value2 = 3
$ws1.Cells.Item(18,4).Value2 = 363.3
$ws1.Cells.Item(18,5).Value2 = 121.47
$ws1.Cells.Item(18,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(18,7).Value2 = "➖ Neutre"

# Row 19: BRVM - SERVICES FINANCIERS
$ws1.Cells.Item(19,1).Value2 = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(19,2).Value2 = 0
$ws1.Cells.Item(19,3).Value2 = 3
$ws1.Cells.Item(19,4).Value2 = 357.04
$ws1.Cells.Item(19,5).Value2 = 119.38
$ws1.Cells.Item(19,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(19,7).Value2 = "➖ Neutre"

# Row 20: BRVM - ENERGIE
$ws1.Cells.Item(20,1).Value2 = "BRVM - ENERGIE"
$ws1.Cells.Item(20,2).Value2 = 0
$ws1.Cells.Item(20,3).Value2 = 3
$ws1.Cells.Item(20,4).Value2 = 323.39
$ws1.Cells.Item(20,5).Value2 = 108.38
$ws1.Cells.Item(20,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(20,7).Value2 = "➖ Neutre"

# Row 21: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Cells.Item(21,1).Value2 = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(21,2).Value2 = 0
$ws1.Cells.Item(21,3).Value2 = 3
$ws1.Cells.Item(21,4).Value2 = 312.58
$ws1.Cells.Item(21,5).Value2 = 104.7
$ws1.Cells.Item(21,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(21,7).Value2 = "➖ Neutre"

# Row 22: BRVM - TELECOMMUNICATIONS
$ws1.Cells.Item(22,1).Value2 = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(22,2).Value2 = 0
$ws1.Cells.Item(22,3).Value2 = 3
$ws1.Cells.Item(22,4).Value2 = 275.74
$ws1.Cells.Item(22,5).Value2 = 91.84
$ws1.Cells.Item(22,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(22,7).Value2 = "➖ Neutre"

# Row 23: SAFCA CI (SAFC)
$ws1.Cells.Item(23,1).Value2 = "SAFCA CI (SAFC)"
$ws1.Cells.Item(23,2).Value2 = 3
$ws1.Cells.Item(23,3).Value2 = 0
$ws1.Cells.Item(23,4).Value2 = 21.16
$ws1.Cells.Item(23,5).Value2 = 7.11
$ws1.Cells.Item(23,6).Value2 = "🟢 Achat"
$ws1.Cells.Item(23,7).Value2 = "✅ Renforcer"

# Row 24: BANK OF AFRICA ML (BOAM)
$ws1.Cells.Item(24,1).Value2 = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(24,2).Value2 = 1
$ws1.Cells.Item(24,3).Value2 = 0
$ws1.Cells.Item(24,4).Value2 = 7.37
$ws1.Cells.Item(24,5).Value2 = 7.37
$ws1.Cells.Item(24,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(24,7).Value2 = "➖ Neutre"

# Row 25: SAPH CI (SPHC)
$ws1.Cells.Item(25,1).Value2 = "SAPH CI (SPHC)"
$ws1.Cells.Item(25,2).Value2 = 1
$ws1.Cells.Item(25,3).Value2 = 0
$ws1.Cells.Item(25,4).Value2 = 4.69
$ws1.Cells.Item(25,5).Value2 = 4.69
$ws1.Cells.Item(25,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(25,7).Value2 = "➖ Neutre"

# Row 26: SICABLE CI (CABC)
$ws1.Cells.Item(26,1).Value2 = "SICABLE CI (CABC)"
$ws1.Cells.Item(26,2).Value2 = 1
$ws1.Cells.Item(26,3).Value2 = 0
$ws1.Cells.Item(26,4).Value2 = 4
$ws1.Cells.Item(26,5).Value2 = 4
$ws1.Cells.Item(26,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(26,7).Value2 = "➖ Neutre"

# Row 27: BERNABE CI (BNBC)
$ws1.Cells.Item(27,1).Value2 = "BERNABE CI (BNBC)"
$ws1.Cells.Item(27,2).Value2 = 1
$ws1.Cells.Item(27,3).Value2 = 1
$ws1.Cells.Item(27,4).Value2 = 3.35
$ws1.Cells.Item(27,5).Value2 = 7.2
$ws1.Cells.Item(27,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(27,7).Value2 = "👀 À surveiller"

# Row 28: VIVO ENERGY CI (SHEC)
$ws1.Cells.Item(28,1).Value2 = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(28,2).Value2 = 1
$ws1.Cells.Item(28,3).Value2 = 0
$ws1.Cells.Item(28,4).Value2 = 3.35
$ws1.Cells.Item(28,5).Value2 = 3.35
$ws1.Cells.Item(28,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(28,7).Value2 = "➖ Neutre"

# Row 29: AIR LIQUIDE CI (SIVC)
$ws1.Cells.Item(29,1).Value2 = "AIR LIQUIDE CI (SIVC)"
$ws1.Cells.Item(29,2).Value2 = 1
$ws1.Cells.Item(29,3).Value2 = 1
$ws1.Cells.Item(29,4).Value2 = 2.94
$ws1.Cells.Item(29,5).Value2 = -1.82
$ws1.Cells.Item(29,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(29,7).Value2 = "👀 À surveiller"

# Row 30: SICOR CI (SICC)
$ws1.Cells.Item(30,1).Value2 = "SICOR CI (SICC)"
$ws1.Cells.Item(30,2).Value2 = 1
$ws1.Cells.Item(30,3).Value2 = 0
$ws1.Cells.Item(30,4).Value2 = 2.8
$ws1.Cells.Item(30,5).Value2 = 2.8
$ws1.Cells.Item(30,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(30,7).Value2 = "➖ Neutre"

# Row 31: BANK OF AFRICA NG (BOAN)
$ws1.Cells.Item(31,1).Value2 = "BANK OF AFRICA NG (BOAN)"
$ws1.Cells.Item(31,2).Value2 = 1
$ws1.Cells.Item(31,3).Value2 = 0
$ws1.Cells.Item(31,4).Value2 = 2.24
$ws1.Cells.Item(31,5).Value2 = 2.24
$ws1.Cells.Item(31,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(31,7).Value2 = "➖ Neutre"

# Row 32: UNIWAX CI (UNXC)
$ws1.Cells.Item(32,1).Value2 = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(32,2).Value2 = 1
$ws1.Cells.Item(32,3).Value2 = 0
$ws1.Cells.Item(32,4).Value2 = 1.75
$ws1.Cells.Item(32,5).Value2 = 1.75
$ws1.Cells.Item(32,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(32,7).Value2 = "➖ Neutre"

# Row 33: NEI-CEDA CI (NEIC)
$ws1.Cells.Item(33,1).Value2 = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(33,2).Value2 = 1
$ws1.Cells.Item(33,3).Value2 = 0
$ws1.Cells.Item(33,4).Value2 = 1.69
$ws1.Cells.Item(33,5).Value2 = 1.69
$ws1.Cells.Item(33,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(33,7).Value2 = "➖ Neutre"

# Row 34: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Cells.Item(34,1).Value2 = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(34,2).Value2 = 1
$ws1.Cells.Item(34,3).Value2 = 1
$ws1.Cells.Item(34,4).Value2 = 0.42
$ws1.Cells.Item(34,5).Value2 = 6.67
$ws1.Cells.Item(34,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(34,7).Value2 = "👀 À surveiller"

# Row 35: TOTAL
$ws1.Cells.Item(35,1).Value2 = "TOTAL"
$ws1.Cells.Item(35,2).Value2 = 0
$ws1.Cells.Item(35,3).Value2 = 3
$ws1.Cells.Item(35,4).Value2 = 0
$ws1.Cells.Item(35,5).Value2 = 0
$ws1.Cells.Item(35,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(35,7).Value2 = "➖ Neutre"

# Row 36: ORANGE COTE D'IVOIRE (ORAC)
$ws1.Cells.Item(36,1).Value2 = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(36,2).Value2 = 0
$ws1.Cells.Item(36,3).Value2 = 1
$ws1.Cells.Item(36,4).Value2 = -1.06
$ws1.Cells.Item(36,5).Value2 = -1.06
$ws1.Cells.Item(36,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(36,7).Value2 = "➖ Neutre"

# Row 37: ECOBANK COTE D''IVOIRE (ECOC)
$ws1.Cells.Item(37,1).Value2 = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(37,2).Value2 = 0
$ws1.Cells.Item(37,3).Value2 = 1
$ws1.Cells.Item(37,4).Value2 = -1.32
$ws1.Cells.Item(37,5).Value2 = -1.32
$ws1.Cells.Item(37,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(37,7).Value2 = "➖ Neutre"

# Row 38: SONATEL SN (SNTS)
$ws1.Cells.Item(38,1).Value2 = "SONATEL SN (SNTS)"
$ws1.Cells.Item(38,2).Value2 = 0
$ws1.Cells.Item(38,3).Value2 = 1
$ws1.Cells.Item(38,4).Value2 = -1.4
$ws1.Cells.Item(38,5).Value2 = -1.4
$ws1.Cells.Item(38,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(38,7).Value2 = "➖ Neutre"

# Row 39: LOTERIE NATIONALE DU BENIN (LNBB)
$ws1.Cells.Item(39,1).Value2 = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(39,2).Value2 = 1
$ws1.Cells.Item(39,3).Value2 = 1
$ws1.Cells.Item(39,4).Value2 = -1.58
$ws1.Cells.Item(39,5).Value2 = 5.88
$ws1.Cells.Item(39,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(39,7).Value2 = "👀 À surveiller"

# Row 40: ONATEL BF (ONTBF)
$ws1.Cells.Item(40,1).Value2 = "ONATEL BF (ONTBF)"
$ws1.Cells.Item(40,2).Value2 = 0
$ws1.Cells.Item(40,3).Value2 = 1
$ws1.Cells.Item(40,4).Value2 = -2.13
$ws1.Cells.Item(40,5).Value2 = -2.13
$ws1.Cells.Item(40,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(40,7).Value2 = "➖ Neutre"

# Row 41: BICI CI (BICC)
$ws1.Cells.Item(41,1).Value2 = "BICI CI (BICC)"
$ws1.Cells.Item(41,2).Value2 = 0
$ws1.Cells.Item(41,3).Value2 = 1
$ws1.Cells.Item(41,4).Value2 = -2.85
$ws1.Cells.Item(41,5).Value2 = -2.85
$ws1.Cells.Item(41,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(41,7).Value2 = "➖ Neutre"

# Row 42: SETAO CI (STAC)
$ws1.Cells.Item(42,1).Value2 = "SETAO CI (STAC)"
$ws1.Cells.Item(42,2).Value2 = 0
$ws1.Cells.Item(42,3).Value2 = 1
$ws1.Cells.Item(42,4).Value2 = -3.45
$ws1.Cells.Item(42,5).Value2 = -3.45
$ws1.Cells.Item(42,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(42,7).Value2 = "➖ Neutre"

# Row 43: ORAGROUP TOGO (ORGT)
$ws1.Cells.Item(43,1).Value2 = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(43,2).Value2 = 0
$ws1.Cells.Item(43,3).Value2 = 2
$ws1.Cells.Item(43,4).Value2 = -5.21
$ws1.Cells.Item(43,5).Value2 = -1.57
$ws1.Cells.Item(43,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(43,7).Value2 = "➖ Neutre"

# Row 44: CIE CI (CIEC)
$ws1.Cells.Item(44,1).Value2 = "CIE CI (CIEC)"
$ws1.Cells.Item(44,2).Value2 = 0
$ws1.Cells.Item(44,3).Value2 = 2
$ws1.Cells.Item(44,4).Value2 = -6.62
$ws1.Cells.Item(44,5).Value2 = -2.71
$ws1.Cells.Item(44,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(44,7).Value2 = "➖ Neutre"

# --- Remove obsolete rows 45-48 (companies dropped from latest dataset) ---
$ws1.Rows("45:48").Delete()

# --- Sheet 2: Top_YTD — refresh YTD progression values (labels unchanged) ---
$ws2.Cells.Item(2,2).Value2 = 424860.98
$ws2.Cells.Item(3,2).Value2 = 41806.25
$ws2.Cells.Item(4,2).Value2 = 39726.94
$ws2.Cells.Item(5,2).Value2 = 33228.73
$ws2.Cells.Item(6,2).Value2 = 31805.6
$ws2.Cells.Item(7,2).Value2 = 30653
$ws2.Cells.Item(8,2).Value2 = 24893.75
$ws2.Cells.Item(9,2).Value2 = 10967
$ws2.Cells.Item(10,2).Value2 = 9829.059999999999
$ws2.Cells.Item(11,2).Value2 = 9013.110000000001

Write-Output "BRVM data refresh applied"